# 004江润洲工作进展更新
# Adds a new weekly progress entry (20240429-20240505) to the
# "004江润洲" worksheet, and leaves the UI selection/active-sheet state
# the way the author left it when saving (004江润洲 tab active, a new
# blank spacer row left on the 001毛彬 tab).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 001毛彬
$ws4 = $wb.Worksheets.Item(4)   # 004江润洲

# --- 004江润洲: append the new week's row (row 5) ---------------------
$ws4.Range("A5").Value = "20240429-20240505"
$ws4.Range("B5").Value = "'12"
$ws4.Range("C5").Value = "1. ASD多站点分类任务`n（1）复现GuidedBackPropogation方法找重要Biomarker`n2. 自我学习`n（1）双周6道题答案整理"
$ws4.Range("C5").WrapText = $true
$ws4.Range("D5").Value = $ws4.Range("D4").Value2
$ws4.Range("E5").Value = $ws4.Range("E4").Value2
$ws4.Range("E5").WrapText = $true
$ws4.Rows.Item(5).RowHeight = 123.75

# --- 001毛彬: leave a blank spacer row (row 17) below the data --------
$ws1.Rows.Item(17).RowHeight = 49.5

# --- restore on-save cursor / selection / active sheet state ----------
$ws1.Activate()
$ws1.Range("A17:XFD17").Select()

$ws4.Activate()
$ws4.Range("D14").Select()
